$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1 -> "Total", D1 -> "Year"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "Year"

# Replace the "FSI Total Normalized" values in column D (rows 2-180) with the year 2024
for ($r = 2; $r -le 180; $r++) {
    $ws.Cells.Item($r, 4).Value = 2024
}

# Remove column E entirely (it held "FSI Total Normalized x4")
$ws.Range("E1:E180").EntireColumn.Delete()
